$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text in B2: remove "/RME" from the third line
$ws.Range("B2").Value = "28% CR+PC/LFM+CDL/H:2`n12% S+SL/LFM+CDL/H:1`n18% S/LFM+CDL/H:1`n15% S/LFBR+CDM/H:1`n27% CR/LFINF+CDL/H:1"

# Wrap text for the updated cell and set a taller row height
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 256

# Match the saved selection state from the authored workbook
[void]$ws.Range("B11").Select()
